$wb = $excel.ActiveWorkbook
$wsData  = $wb.Worksheets.Item("Table Data")
$wsLogin = $wb.Worksheets.Item("Login")

# --- "Table Data" sheet -----------------------------------------------
# Header rename: "Resident (Static)" -> "Resident"
$wsData.Range("F1").Value = "Resident"

# Populate the (previously empty) "Resident" column with the resident
# name tied to each transfer-in row.
$wsData.Range("F2").Value  = "Arvind Nath"
$wsData.Range("F3").Value  = "Windsor Charles"
$wsData.Range("F4").Value  = "Wilbur Smith"
$wsData.Range("F5").Value  = "Jerry RAC004"
$wsData.Range("F6").Value  = "Ita Rooney"
$wsData.Range("F7").Value  = "Robyn Dhar"
$wsData.Range("F8").Value  = "Robert Jones"
$wsData.Range("F9").Value  = "Ted Bryan"
$wsData.Range("F10").Value = "Ita Rooney"
$wsData.Range("F11").Value = "Michael ILU"
$wsData.Range("F12").Value = "Pamela Butler"
$wsData.Range("F13").Value = "Perry Grant"
$wsData.Range("F14").Value = "David Springer"
$wsData.Range("F15").Value = "Robyn Dhar"
$wsData.Range("F16").Value = "Matilda Kerr"

# One row's resident name picked up a distinct (re-applied) font.
$wsData.Range("F13").Font.Name = "Arial"

# Column F/G used to be a pair of narrow spacer columns (width 1.39); now
# that F holds real data it (and the next column G) become normal, wider
# data columns.
$wsData.Columns.Item(6).ColumnWidth = 17.636666666666667
$wsData.Columns.Item(7).ColumnWidth = 19.17666666666667

# --- "Login" sheet ------------------------------------------------------
# Move the cursor/selection, then switch the active sheet back to
# "Table Data" with its own new selection -- matches the saved view state.
$wsLogin.Range("E27").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("F2").Select() | Out-Null
